$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "mngr549847"
$ws.Range("B1").Value = "surUpeg"

$ws.Range("D8").Select()
